$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 357 (pushing the
# existing rows 357-374 down to 359-376). Row insertion in this sheet
# inherits the column-D date style ("s=2") from the surrounding rows, so
# the new rows already pick up the correct number format for the Fecha
# column.
$ws.Rows.Item(357).Insert()
$ws.Rows.Item(357).Insert()

# New row 357 - weekly Cilantro price entry (Provincia de Cautín)
$ws.Cells.Item(357, 1).Value = 10
$ws.Cells.Item(357, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(357, 3).Value = "La Araucanía"
$ws.Cells.Item(357, 4).Value = 44753
$ws.Cells.Item(357, 5).Value = 9
$ws.Cells.Item(357, 6).Value = 100112040
$ws.Cells.Item(357, 7).Value = "Cilantro"
$ws.Cells.Item(357, 8).Value = "Sin especificar"
$ws.Cells.Item(357, 9).Value = "Primera"
$ws.Cells.Item(357, 10).Value = 50
$ws.Cells.Item(357, 11).Value = 6000
$ws.Cells.Item(357, 12).Value = 6000
$ws.Cells.Item(357, 13).Value = 6000
$ws.Cells.Item(357, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(357, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(357, 16).Value = 3000
$ws.Cells.Item(357, 17).Value = 2
$ws.Cells.Item(357, 18).Value = "Hortaliza"

# New row 358 - weekly Cilantro price entry (Región Metropolitana)
$ws.Cells.Item(358, 1).Value = 10
$ws.Cells.Item(358, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(358, 3).Value = "La Araucanía"
$ws.Cells.Item(358, 4).Value = 44753
$ws.Cells.Item(358, 5).Value = 9
$ws.Cells.Item(358, 6).Value = 100112040
$ws.Cells.Item(358, 7).Value = "Cilantro"
$ws.Cells.Item(358, 8).Value = "Sin especificar"
$ws.Cells.Item(358, 9).Value = "Primera"
$ws.Cells.Item(358, 10).Value = 60
$ws.Cells.Item(358, 11).Value = 4300
$ws.Cells.Item(358, 12).Value = 4300
$ws.Cells.Item(358, 13).Value = 4300
$ws.Cells.Item(358, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(358, 15).Value = "Región Metropolitana"
$ws.Cells.Item(358, 16).Value = 2150
$ws.Cells.Item(358, 17).Value = 2
$ws.Cells.Item(358, 18).Value = "Hortaliza"
